$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1871657754010695
$ws.Range("C2").Value = 0.5668449197860963
$ws.Range("O2").Value = 0.0053475935828877
$ws.Range("P2").Value = 0.1176470588235294
$ws.Range("S2").Value = 0.1229946524064171
$ws.Range("B3").Value = 0.008849557522123894
$ws.Range("C3").Value = 0.06194690265486726
$ws.Range("J3").Value = 0.008849557522123894
$ws.Range("P3").Value = 0.7522123893805309
$ws.Range("S3").Value = 0.168141592920354
$ws.Range("P4").Value = 0.6333333333333333
$ws.Range("S4").Value = 0.3666666666666666
$ws.Range("B6").Value = 0.0427807486631016
$ws.Range("D6").Value = 0.0053475935828877
$ws.Range("F6").Value = 0.06951871657754011
$ws.Range("J6").Value = 0.2032085561497326
$ws.Range("O6").Value = 0.03208556149732621
$ws.Range("Q6").Value = 0.2245989304812834
$ws.Range("R6").Value = 0.0855614973262032
$ws.Range("S6").Value = 0.3368983957219251
$ws.Range("B7").Value = 0.101123595505618
$ws.Range("D7").Value = 0.02247191011235955
$ws.Range("E7").Value = 0.005617977528089887
$ws.Range("F7").Value = 0.0449438202247191
$ws.Range("J7").Value = 0.1067415730337079
$ws.Range("O7").Value = 0.05056179775280899
$ws.Range("Q7").Value = 0.151685393258427
$ws.Range("R7").Value = 0.07865168539325842
$ws.Range("S7").Value = 0.4382022471910113
$ws.Range("B8").Value = 0.07242339832869081
$ws.Range("D8").Value = 0.01392757660167131
$ws.Range("F8").Value = 0.05292479108635097
$ws.Range("J8").Value = 0.08356545961002786
$ws.Range("O8").Value = 0.03621169916434541
$ws.Range("Q8").Value = 0.1949860724233983
$ws.Range("R8").Value = 0.0947075208913649
$ws.Range("S8").Value = 0.4512534818941504
$ws.Range("B9").Value = 0.07251908396946564
$ws.Range("D9").Value = 0.01145038167938931
$ws.Range("F9").Value = 0.08015267175572519
$ws.Range("J9").Value = 0.1106870229007634
$ws.Range("O9").Value = 0.03053435114503817
$ws.Range("Q9").Value = 0.1946564885496183
$ws.Range("R9").Value = 0.1297709923664122
$ws.Range("S9").Value = 0.3702290076335878
$ws.Range("B10").Value = 0.07738095238095238
$ws.Range("D10").Value = 0.01785714285714286
$ws.Range("E10").Value = 0.000992063492063492
$ws.Range("F10").Value = 0.07738095238095238
$ws.Range("J10").Value = 0.1041666666666667
$ws.Range("O10").Value = 0.02083333333333333
$ws.Range("Q10").Value = 0.2023809523809524
$ws.Range("R10").Value = 0.07539682539682539
$ws.Range("S10").Value = 0.4236111111111111
$ws.Range("G11").Value = 0.1266968325791855
$ws.Range("J11").Value = 0.05429864253393665
$ws.Range("K11").Value = 0.1447963800904978
$ws.Range("L11").Value = 0.669683257918552
$ws.Range("S11").Value = 0.004524886877828055
$ws.Range("G12").Value = 0.8092105263157895
$ws.Range("J12").Value = 0.131578947368421
$ws.Range("K12").Value = 0.01973684210526316
$ws.Range("L12").Value = 0.03947368421052631
$ws.Range("G13").Value = 0.7777777777777778
$ws.Range("J13").Value = 0.1555555555555556
$ws.Range("S13").Value = 0.06666666666666667
$ws.Range("H15").Value = 0.1822033898305085
$ws.Range("I15").Value = 0.06779661016949153
$ws.Range("J15").Value = 0.3474576271186441
$ws.Range("K15").Value = 0.05508474576271186
$ws.Range("M15").Value = 0.008474576271186441
$ws.Range("O15").Value = 0.05932203389830509
$ws.Range("S15").Value = 0.2796610169491525
$ws.Range("F16").Value = 0.03278688524590164
$ws.Range("H16").Value = 0.1311475409836066
$ws.Range("I16").Value = 0.1885245901639344
$ws.Range("J16").Value = 0.3278688524590164
$ws.Range("K16").Value = 0.07377049180327869
$ws.Range("M16").Value = 0.04098360655737705
$ws.Range("O16").Value = 0.1147540983606557
$ws.Range("S16").Value = 0.09016393442622951
$ws.Range("F17").Value = 0.01518987341772152
$ws.Range("H17").Value = 0.179746835443038
$ws.Range("I17").Value = 0.1468354430379747
$ws.Range("J17").Value = 0.3949367088607595
$ws.Range("K17").Value = 0.05063291139240506
$ws.Range("M17").Value = 0.03291139240506329
$ws.Range("O17").Value = 0.08354430379746836
$ws.Range("S17").Value = 0.09620253164556962
$ws.Range("F18").Value = 0.02272727272727273
$ws.Range("H18").Value = 0.1477272727272727
$ws.Range("I18").Value = 0.1534090909090909
$ws.Range("J18").Value = 0.4375
$ws.Range("K18").Value = 0.08522727272727272
$ws.Range("M18").Value = 0.01136363636363636
$ws.Range("O18").Value = 0.07954545454545454
$ws.Range("S18").Value = 0.0625
$ws.Range("F19").Value = 0.009700176366843033
$ws.Range("H19").Value = 0.1825396825396825
$ws.Range("I19").Value = 0.1208112874779541
$ws.Range("J19").Value = 0.3677248677248677
$ws.Range("K19").Value = 0.109347442680776
$ws.Range("M19").Value = 0.02028218694885361
$ws.Range("N19").Value = 0.0008818342151675485
$ws.Range("O19").Value = 0.06966490299823633
$ws.Range("S19").Value = 0.119047619047619
